$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Active Cases"

$data = @(
    @('3398 BlueCross Elly Kay Mordialloc', 34),
    @('3601 Baptcare Westhaven community', 14),
    @('3646 Mornington Bay Care Community Mount Martha', 17),
    @('3653 Fronditha Thalpori St Albans Aged Care', 20),
    @('3975 Aurrum Aged Care Brunswick West', 10),
    @('4257 BlueCross The Gables Camberwell', 18),
    @('4295 Hope Aged Care Sunshine West', 21),
    @('4314 Estia Health Ardeer', 18),
    @('44304 Brighton Primary School Brighton', 12),
    @('44380 Plenty Parklands Primary School Mill Park', 11),
    @('44414 Buangor Primary School Buangor', 20),
    @('44584 Badger Creek Primary School Badger Creek', 14),
    @('44666 Gardenvale Primary School Senior School Campus Brighton East', 15),
    @('44761 Coburg North Primary School Coburg', 24),
    @('44828 Cheltenham East Primary School Cheltenham', 10),
    @('44930 St Albans Heights Primary School', 14),
    @('44950 Templestowe Valley Primary School Templestowe Lower', 15),
    @('44979 Campbellfield Heights Primary School Campbellfield', 21),
    @('45248 Brookside P-9 College Caroline Springs', 15),
    @('45257 Roxburgh Rise Primary School Roxburgh Park', 10),
    @('45265 Cranbourne Carlisle Primary School', 10),
    @('45267 Epping Views Primary School Epping', 13),
    @('45648 St Brendans Primary School Shepparton', 10),
    @('45677 St. Peter''s Primary School Bentleigh East', 11),
    @('45708 St Augustine''s Primary School Yarraville', 12),
    @('45719 St Joseph''s Primary School Numurkah', 10),
    @('4574 Village Glen Aged Care Residences Mornington', 11),
    @('45797 St John''s Primary School Footscray', 10),
    @('45812 St Mary''s Primary School Hampton', 10),
    @('45902 Mother of God Primary School Ardeer', 11),
    @('45988 St Macartan''s Parish Primary School', 10),
    @('46052 St. Francis of Assisi Primary School Mill Park', 14),
    @('46105 Christ the Priest Primary School Caroline Springs', 12),
    @('46115 St Luke''s Catholic Primary School Shepparton North', 17),
    @('46221 Bialik College Hawthorn', 11),
    @('46306 King''s College Warrnambool', 10),
    @('46328 Ilim College Inverloch Crescent Dallas', 11),
    @('50722 Chairo Christian School Leongatha Campus Leongatha', 11),
    @('Kororoit Creek Primary School Burnside Heights', 20),
    @('Mambourin Enterprises Allara Deer Park', 19),
    @('North St Kilda Children''s Centre St Kilda', 11),
    @('Oakleigh South Primary School Oakleigh South', 11),
    @('Rosebud Primary School Rosebud', 15),
    @('Springside Primary School Caroline Springs', 17),
    @('St Louis de Montfort''s School Aspendale', 12),
    @('The Village Early Learning Centre Sandringham', 17),
    @('Torquay Hotel Torquay', 28)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
